$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Rename the values in row 1 (B1, C1, D1): "Value 1" -> "Data1", "Value 2" -> "Data2", "Value 3" -> "Data3"
$ws.Range("B1").Value = "Data1"
$ws.Range("C1").Value = "Data2"
$ws.Range("D1").Value = "Data3"

# Update the active cell selection on the TestData sheet to D1
$ws.Activate()
$ws.Range("D1").Select()
